$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Objective" column header (shared string #10)
$ws.Range("T4").Value = "Objective"

# New formulas in column T: Objective = R + S (extreme-point domain check)
$ws.Range("T5:T8").FormulaR1C1 = "=RC[-2]+RC[-1]"

# Emphasize the final objective-domain value
$ws.Range("T8").Font.Bold = $true

# Match the new selection left behind by the editing session
$ws.Range("T8").Select()

# Page orientation was (re)applied during the edit session
$ws.PageSetup.Orientation = 1
